# Merge the title's two runs ("PIWG" + " Action Item XX-XXX") into a single
# run reading "PIWG Action Item XX-XXX", and drop the stale spell-check
# "err" flag that was stuck on the first run, while keeping a clean
# "dirty=0" marker (matching the target OOXML: <a:rPr lang="en-US" dirty="0"/>).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)          # "Title 1"
$tr = $sh.TextFrame.TextRange

# Drop the "PIWG" run's inherited err="1" flag by deleting it and typing it
# back in front of the remaining text, which inherits the clean formatting
# of the run that follows it.
$tr.Characters(1, 4).Text = ""
$null = $tr.InsertBefore("PIWG")

# Re-assert the final wording and collapse everything to one run with
# uniform (clean) formatting. Routing through an unrelated placeholder
# value (sharing no prefix with the final text) forces a full rebuild of
# the paragraph's runs instead of an incremental "append a new run" edit.
$tr.Text = "TEMP_MERGE_VALUE"
$tr.Text = "PIWG Action Item XX-XXX"
